# Natmi following Dr Hou advice
# Updates the Gdnf-Gfra1 LR-pair sheet:
#  - Row 2 (sCs -> ECs) values recomputed
#  - Row 3 (sCs -> FAPs) is new content (used to be the only "other" row, now re-targeted to FAPs)
#  - Row 4 (sCs -> sCs) is an entirely new row

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: Target cluster "ECs" ---
$ws.Cells.Item(2, 4).Value = "ECs"
$ws.Cells.Item(2, 5).Value = 3
$ws.Cells.Item(2, 6).Value = 1
$ws.Cells.Item(2, 7).Value = 0.697917
$ws.Cells.Item(2, 8).Value = 2.093751
$ws.Cells.Item(2, 9).Value = 1
$ws.Cells.Item(2, 10).Value = 1
$ws.Cells.Item(2, 11).Value = 2
$ws.Cells.Item(2, 12).Value = 0.6666666666666666
$ws.Cells.Item(2, 13).Value = 0.06547366666666667
$ws.Cells.Item(2, 14).Value = 0.196421
$ws.Cells.Item(2, 15).Value = 0.002125877360986814
$ws.Cells.Item(2, 16).Value = 0.002125877360986814
$ws.Cells.Item(2, 17).Value = 0.045695185019
$ws.Cells.Item(2, 18).Value = 0.4112566651710001
$ws.Cells.Item(2, 19).Value = 0.002125877360986814
$ws.Cells.Item(2, 20).Value = 0.002125877360986814

# --- Row 3: Target cluster "FAPs" (re-purposed from the old single extra row) ---
$ws.Cells.Item(3, 1).Value = "sCs"
$ws.Cells.Item(3, 2).Value = "Gdnf"
$ws.Cells.Item(3, 3).Value = "Gfra1"
$ws.Cells.Item(3, 4).Value = "FAPs"
$ws.Cells.Item(3, 5).Value = 3
$ws.Cells.Item(3, 6).Value = 1
$ws.Cells.Item(3, 7).Value = 0.697917
$ws.Cells.Item(3, 8).Value = 2.093751
$ws.Cells.Item(3, 9).Value = 1
$ws.Cells.Item(3, 10).Value = 1
$ws.Cells.Item(3, 11).Value = 3
$ws.Cells.Item(3, 12).Value = 1
$ws.Cells.Item(3, 13).Value = 25.94643066666667
$ws.Cells.Item(3, 14).Value = 77.839292
$ws.Cells.Item(3, 15).Value = 0.8424597607080814
$ws.Cells.Item(3, 16).Value = 0.8424597607080814
$ws.Cells.Item(3, 17).Value = 18.108455051588
$ws.Cells.Item(3, 18).Value = 162.976095464292
$ws.Cells.Item(3, 19).Value = 0.8424597607080814
$ws.Cells.Item(3, 20).Value = 0.8424597607080814

# --- Row 4: Target cluster "sCs" (new row) ---
$ws.Cells.Item(4, 1).Value = "sCs"
$ws.Cells.Item(4, 2).Value = "Gdnf"
$ws.Cells.Item(4, 3).Value = "Gfra1"
$ws.Cells.Item(4, 4).Value = "sCs"
$ws.Cells.Item(4, 5).Value = 3
$ws.Cells.Item(4, 6).Value = 1
$ws.Cells.Item(4, 7).Value = 0.697917
$ws.Cells.Item(4, 8).Value = 2.093751
$ws.Cells.Item(4, 9).Value = 1
$ws.Cells.Item(4, 10).Value = 1
$ws.Cells.Item(4, 11).Value = 3
$ws.Cells.Item(4, 12).Value = 1
$ws.Cells.Item(4, 13).Value = 4.786517
$ws.Cells.Item(4, 14).Value = 14.359551
$ws.Cells.Item(4, 15).Value = 0.1554143619309319
$ws.Cells.Item(4, 16).Value = 0.1554143619309319
$ws.Cells.Item(4, 17).Value = 3.340591585089
$ws.Cells.Item(4, 18).Value = 30.065324265801
$ws.Cells.Item(4, 19).Value = 0.1554143619309319
$ws.Cells.Item(4, 20).Value = 0.1554143619309319
